$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212)
$newDate = (Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
